# Update "want-to-go" counts (column F) and a couple of min-price values
# (column G) across the four sheets, matching the refreshed data snapshot
# from the site generator (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 284
$ws1.Range("F4").Value = 638
$ws1.Range("F5").Value = 2829
$ws1.Range("F8").Value = 6561
$ws1.Range("F11").Value = 5102
$ws1.Range("F12").Value = 10
$ws1.Range("F14").Value = 2706
$ws1.Range("F17").Value = 1247
$ws1.Range("G17").Value = 68
$ws1.Range("F21").Value = 1116
$ws1.Range("F23").Value = 1406
$ws1.Range("F24").Value = 1070
$ws1.Range("F25").Value = 2129
$ws1.Range("F26").Value = 1352
$ws1.Range("F28").Value = 57
$ws1.Range("F29").Value = 1016
$ws1.Range("F30").Value = 51
$ws1.Range("F31").Value = 126
$ws1.Range("F32").Value = 1549
$ws1.Range("F33").Value = 12
$ws1.Range("F35").Value = 1625
$ws1.Range("F36").Value = 1105
$ws1.Range("F39").Value = 327
$ws1.Range("F40").Value = 2342
$ws1.Range("F41").Value = 2604
$ws1.Range("F48").Value = 120

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 21
$ws2.Range("F8").Value = 335
$ws2.Range("F10").Value = 172
$ws2.Range("F11").Value = 104
$ws2.Range("F26").Value = 432

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 525
$ws3.Range("F6").Value = 1711
$ws3.Range("F8").Value = 1562
$ws3.Range("F10").Value = 2603
$ws3.Range("F11").Value = 909
$ws3.Range("F12").Value = 797

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 525
$ws4.Range("F4").Value = 1711
$ws4.Range("F5").Value = 638
$ws4.Range("F6").Value = 2829
$ws4.Range("F7").Value = 1562
$ws4.Range("F9").Value = 2603
$ws4.Range("F10").Value = 6562
$ws4.Range("F11").Value = 909
$ws4.Range("F12").Value = 797
$ws4.Range("F13").Value = 5102
$ws4.Range("F14").Value = 2706
$ws4.Range("F17").Value = 1247
$ws4.Range("G17").Value = 68
$ws4.Range("F20").Value = 335
$ws4.Range("F21").Value = 1116
$ws4.Range("F23").Value = 104
$ws4.Range("F25").Value = 1406
$ws4.Range("F26").Value = 1070
$ws4.Range("F27").Value = 2129
$ws4.Range("F28").Value = 1352
$ws4.Range("F30").Value = 57
$ws4.Range("F32").Value = 1016
$ws4.Range("F33").Value = 51
$ws4.Range("F35").Value = 1549
$ws4.Range("F37").Value = 1105
$ws4.Range("F39").Value = 432
$ws4.Range("F40").Value = 327
$ws4.Range("F43").Value = 2342
$ws4.Range("F44").Value = 2604
